$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    3 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    5 = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.242251378316819)
    6 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    7 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    8 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    9 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
